$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare row 131 formatting (copy from row 130 style: bold/bordered id col, date col) ---
$ws.Range("A130").Copy()
$ws.Range("A131").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E130").Copy()
$ws.Range("E131").PasteSpecial(-4122)   # xlPasteFormats

# --- Update existing row 130 (id 128): Hyderabad FC vs Kerala Blasters, 2024-04-12 ---
$ws.Range("B130").Value = 7749472
$ws.Range("E130").Value = 45394.45833333334
$ws.Range("F130").Value = "Hyderabad FC"
$ws.Range("G130").Value = "Kerala Blasters"
$ws.Range("H130").Value = 1
$ws.Range("I130").Value = 3
$ws.Range("J130").Value = "A"
$ws.Range("K130").Value = 5.25
$ws.Range("L130").Value = 4
$ws.Range("M130").Value = 1.6
$ws.Range("N130").Value = 3.3
$ws.Range("O130").Value = 3.6
$ws.Range("P130").Value = 2.05
$ws.Range("Q130").Value = 0.25
$ws.Range("R130").Value = 2.025
$ws.Range("S130").Value = 1.825
$ws.Range("T130").Value = 3
$ws.Range("U130").Value = 2
$ws.Range("V130").Value = 1.85
$ws.Range("W130").Value = -1
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 1.05
$ws.Range("Z130").Value = -1
$ws.Range("AA130").Value = 0.825
$ws.Range("AB130").Value = 1
$ws.Range("AC130").Value = -1

# --- New row 131 (id 129): Mohun Bagan SG vs Mumbai City FC, 2024-04-15 ---
$ws.Range("A131").Value = 129
$ws.Range("B131").Value = 7749764
$ws.Range("C131").Value = "India Super League"
$ws.Range("D131").Value = "India Super League"
$ws.Range("E131").Value = 45397.45833333334
$ws.Range("F131").Value = "Mohun Bagan SG"
$ws.Range("G131").Value = "Mumbai City FC"
$ws.Range("K131").Value = 2.5
$ws.Range("L131").Value = 3.4
$ws.Range("M131").Value = 2.7
$ws.Range("N131").Value = 2.5
$ws.Range("O131").Value = 3.4
$ws.Range("P131").Value = 2.7
$ws.Range("Q131").Value = 0
$ws.Range("R131").Value = 1.825
$ws.Range("S131").Value = 1.975
$ws.Range("T131").Value = 2.75
$ws.Range("U131").Value = 1.975
$ws.Range("V131").Value = 1.825
$ws.Range("W131").Value = 0
$ws.Range("X131").Value = 0
$ws.Range("Y131").Value = 0
$ws.Range("Z131").Value = 0
$ws.Range("AA131").Value = 0
